$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently follows the H1
#    title ("Play Chilli Pop for Free: Exciting Mexican-Themed Slot Game").
# ---------------------------------------------------------------------------
$metaFound = $d.Content.Find.Execute("Meta description", $false, $false,
    $false, $false, $false, $true, 1, $false, "", 0)
if ($metaFound -and $d.Content.Find.Found) {
    $metaRange = $d.Content
    $metaPara = $metaRange.Paragraphs.Item(1)
    $metaPara.Range.Delete()
} else {
    # Fallback: the meta description paragraph is the 2nd paragraph in the doc
    $d.Paragraphs.Item(2).Range.Delete()
}

# ---------------------------------------------------------------------------
# 2. Split the final paragraph (the italic AI image-prompt paragraph) into
#    two paragraphs: a new bold "Play Chilli Pop for Free..." paragraph,
#    followed by the existing (still-italic) paragraph.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertAt = $lastPara.Range.Start
$insertionRange = $d.Range($insertAt, $insertAt)

$newParaXml = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Chilli Pop for Free: Exciting Mexican-Themed Slot Game</w:t></w:r></w:p>' +
    '<w:p><w:bookmarkStart w:id="998877" w:name="zzzTempSplitMarker"/><w:bookmarkEnd w:id="998877"/></w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($newParaXml)

# Clean up the temporary bookmark used purely to anchor the paragraph split.
$tempBookmark = $d.Bookmarks.Item("zzzTempSplitMarker")
$tempBookmark.Delete()

# ---------------------------------------------------------------------------
# 3. Replace the old AI image-prompt text (now in the last paragraph) with
#    the new meta-description sentence, keeping its italic formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count)
$finalPara.Range.Find.Execute(
    "Create a vibrant feature image for Chilli Pop that showcases the game's Mexican theme and fun characters. The image should be in a cartoon style and feature a happy Maya warrior with glasses as the main focus. Surround the warrior with symbols from the game, such as tomatoes, garlic, peppers, and onions that have been transformed into wacky characters. Use bright colors to make the image pop and ensure that it captures the spirit of the game's exciting gameplay and cluster-based winning combinations. Add the Chilli Pop logo to the center of the image to tie it all together and make it clear which game it represents.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Chilli Pop, a colorful Mexican cuisine-inspired online slot game. Play for free and enjoy free spins, bonus features, and multipliers.",
    2)
